# Apply the updated crypto price/volume snapshot to Sheet1.
# Cells in column D whose text looks like a plain number are written with a
# leading "'" (Excel quote-prefix) so COM stores them as text, matching the
# original inline/shared-string cells instead of being coerced to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.765.33'
$ws.Range("E2").Value = '  +0.35%  '

$ws.Range("D3").Value = '1.644.76'
$ws.Range("E3").Value = '  +0.09%  '

$ws.Range("E4").Value = '  +0.39%  '

$ws.Range("D5").Value = "'216.63"
$ws.Range("E5").Value = '  +0.73%  '

$ws.Range("D6").Value = "'0.502"
$ws.Range("E6").Value = '  -0.69%  '

$ws.Range("E7").Value = '  +0.45%  '

$ws.Range("E8").Value = '  -0.28%  '

$ws.Range("E9").Value = '  -0.11%  '

$ws.Range("E10").Value = '  -0.84%  '

$ws.Range("E11").Value = '  +0.23%  '

$ws.Range("D12").Value = '1.656.42'
$ws.Range("E12").Value = '  +0.63%  '

$ws.Range("E13").Value = '  -0.93%  '

$ws.Range("E14").Value = '  -0.63%  '

$ws.Range("D15").Value = "'64.66"
$ws.Range("E15").Value = '  -1.31%  '

$ws.Range("D16").Value = '26.771.22'
$ws.Range("E16").Value = '  +0.20%  '

$ws.Range("D17").Value = '0.0₃0736'
$ws.Range("E17").Value = '  -1.75%  '

$ws.Range("D18").Value = "'214.15"
$ws.Range("E18").Value = '  -1.19%  '

$ws.Range("E19").Value = '  +0.42%  '

$ws.Range("E20").Value = '  +0.70%  '

$ws.Range("E21").Value = '  +12.51%  '

$ws.Range("D22").Value = "'6.25"
$ws.Range("E22").Value = '  -1.02%  '

$ws.Range("E23").Value = '  -2.15%  '

$ws.Range("D24").Value = "'146.12"

$ws.Range("E25").Value = '  +0.58%  '

$ws.Range("D27").Value = "'7.14"
$ws.Range("E27").Value = '  -0.50%  '

$ws.Range("E28").Value = '  -0.98%  '

$ws.Range("E29").Value = '  -1.63%  '

$ws.Range("E30").Value = '  +0.95%  '

$ws.Range("E31").Value = '  -0.41%  '

$ws.Range("D32").Value = "'3.00"
$ws.Range("E32").Value = '  -1.52%  '

$ws.Range("D33").Value = '1.288.35'
$ws.Range("E33").Value = '  +1.14%  '

$ws.Range("E34").Value = '  -0.64%  '

$ws.Range("E36").Value = '  -3.05%  '

$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("E38").Value = '  -1.27%  '

$ws.Range("D40").Value = "'0.805"
$ws.Range("E40").Value = '  -1.21%  '

$ws.Range("E41").Value = '  -1.01%  '

$ws.Range("E42").Value = '  -2.50%  '

$ws.Range("D43").Value = '1.789.85'
$ws.Range("E43").Value = '  +0.41%  '

$ws.Range("D44").Value = "'61.60"
$ws.Range("E44").Value = '  +2.69%  '

$ws.Range("E45").Value = '  -0.82%  '

$ws.Range("E46").Value = '  +0.20%  '

$ws.Range("D47").Value = '0.0₆0103'
$ws.Range("E47").Value = '  -1.56%  '

$ws.Range("D48").Value = "'0.0521"
$ws.Range("E48").Value = '  +1.01%  '

$ws.Range("D49").Value = "'7.67"
$ws.Range("E49").Value = '  -1.79%  '

$ws.Range("D50").Value = "'0.0970"
$ws.Range("E50").Value = '  -0.23%  '

$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = "'1.01"
$ws.Range("E51").Value = '  +0.64%  '
